# Updates on Sprint documents
# Replace the sample "User Story" rows with the real Agile user stories for
# the "Aposter" persona, add the 4th story that used to be blank, resize the
# bottom screenshot picture, and refresh the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Row 6 : Story 1 --------------------------------------------------
# (write C/D before B so new shared-string entries land in the same order
#  the original author produced them in)
$ws.Range("C6").Value = "View a 2 lines report of game"
$ws.Range("D6").Value = "have useful stats"
$ws.Range("B6").Value = "Aposter"

# -- Row 7 : Story 2 --------------------------------------------------
$ws.Range("B7").Value = "Aposter"
$ws.Range("C7").Value = "View a series report ordered by day"
$ws.Range("D7").Value = "have important info easily on hand"

# -- Row 8 : Story 3 --------------------------------------------------
$ws.Range("B8").Value = "Aposter"
$ws.Range("C8").Value = "Order the series by sheets depending of the start day"
$ws.Range("D8").Value = "have more order on the today and upcoming serie games"

# -- Row 9 : Story 4 (previously blank) -------------------------------
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Aposter"
$ws.Range("C9").Value = "Have info from MLB and Covers on my DB"
$ws.Range("D9").Value = "access easily important data"

# -- Shrink the bottom screenshot picture by 51pt (~0.7in) ------------
$sh = $ws.Shapes.Item(2)
$sh.Width = $sh.Width - 51.0

# -- Move the active selection to C9, matching the saved view ---------
[void]$ws.Range("C9").Select()
